$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E values (shifted/rotated versions of column C, plus a brand-new label)
$ws.Range("E1").Value = 43
$ws.Range("E4").Value = " State Diagram for Pedestrian And Car TLS"
$ws.Range("E5").Value = "Coding for Arduino Circuit Pedestrian &Car traffic light"
$ws.Range("E6").Value = "Class Diagram for Pedestrian And Car TLS"
$ws.Range("E7").Value = "Coding and connecting circuit through TinkerCAD"

# Widen column C to match new layout
$ws.Columns("C").ColumnWidth = 44.3

# Update selection to match target file
$ws.Range("D4").Select()
